$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:2").Insert()

$ws.Range("A2").Value = "(пайыз менен)"
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"
